# TC01_C3DC_phs001437_SexAtBirth-Unknown.xlsx
# "Updated C3DC phs002371 input query"
#
# The "TreatmentTab" row's query (cell B5 on Sheet1) needs an extra filter
# clause added to its WHERE predicate so that only participants with a
# non-null treatment id are returned:
#
#   ... AND prt.sex_at_birth = 'Unknown'
#   ORDER BY ...
#
# becomes
#
#   ... AND prt.sex_at_birth = 'Unknown' AND trt.treatment_id IS NOT NULL
#   ORDER BY ...

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cell = $ws.Range("B5")
$query = $cell.Value2

$find    = "std.dbgap_accession = 'phs001437' AND prt.sex_at_birth = 'Unknown'`nORDER BY"
$replace = "std.dbgap_accession = 'phs001437' AND prt.sex_at_birth = 'Unknown' AND trt.treatment_id IS NOT NULL`nORDER BY"

if ($query.Contains($find)) {
    $cell.Value2 = $query.Replace($find, $replace)
} else {
    Write-Host "WARNING: expected WHERE clause not found in B5 - leaving untouched"
}

Write-Host $ws.Range("B5").Value2
